# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed May 31 04:32:59 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force text storage for values that would otherwise be
# auto-coerced to a Number by Excel (e.g. '1.000' -> 1), mirroring
# the quote-prefix a user gets by typing a leading apostrophe.
function Set-CellText($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# Row 2
$ws.Range("D2").Value = '27.667.85'
$ws.Range("E2").Value = '  -0.67%  '

# Row 3
$ws.Range("D3").Value = '1.895.95'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
Set-CellText "D4" '1.001'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
Set-CellText "D5" '310.33'
$ws.Range("E5").Value = '  -0.81%  '

# Row 6
$ws.Range("E6").Value = '  -0.14%  '

# Row 7
Set-CellText "D7" '0.5242'
$ws.Range("E7").Value = '  +4.33%  '

# Row 8
Set-CellText "D8" '0.3805'
$ws.Range("E8").Value = '  -0.20%  '

# Row 9
Set-CellText "D9" '0.07235'
$ws.Range("E9").Value = '  -0.71%  '

# Row 10
Set-CellText "D10" '21.11'
$ws.Range("E10").Value = '  +1.14%  '

# Row 11
Set-CellText "D11" '0.9019'
$ws.Range("E11").Value = '  -0.98%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.892.27'
$ws.Range("E12").Value = '  -0.29%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText "D13" '0.07634'
$ws.Range("E13").Value = '  -0.19%  '

# Row 14
Set-CellText "D14" '5.428'
$ws.Range("E14").Value = '  -1.03%  '

# Row 15
Set-CellText "D15" '91.62'
$ws.Range("E15").Value = '  +0.30%  '

# Row 16
Set-CellText "D16" '1.001'
$ws.Range("E16").Value = '  -0.14%  '

# Row 17
Set-CellText "D17" '0.000008670'
$ws.Range("E17").Value = '  -0.54%  '

# Row 18
Set-CellText "D18" '14.33'
$ws.Range("E18").Value = '  -1.31%  '

# Row 19
Set-CellText "D19" '1.000'
$ws.Range("E19").Value = '  -0.17%  '

# Row 20
$ws.Range("D20").Value = '27.701.00'
$ws.Range("E20").Value = '  -0.65%  '

# Row 21
Set-CellText "D21" '5.154'
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").Value = '2.115.87'
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("E23").Value = '  -0.24%  '

# Row 24
Set-CellText "D24" '6.590'
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
Set-CellText "D25" '153.08'
$ws.Range("E25").Value = '  -0.92%  '

# Row 26
Set-CellText "D26" '1.865'
$ws.Range("E26").Value = '  +0.43%  '

# Row 27
Set-CellText "D27" '18.25'
$ws.Range("E27").Value = '  -0.68%  '

# Row 28
Set-CellText "D28" '2.186'

# Row 29
Set-CellText "D29" '114.12'
$ws.Range("E29").Value = '  -0.96%  '

# Row 30
Set-CellText "D30" '4.837'
$ws.Range("E30").Value = '  -1.98%  '

# Row 31
Set-CellText "D31" '4.800'
$ws.Range("E31").Value = '  +3.38%  '

# Row 32
Set-CellText "D32" '0.09130'
$ws.Range("E32").Value = '  +1.78%  '

# Row 33
Set-CellText "D33" '0.05269'
$ws.Range("E33").Value = '  -0.11%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText "D34" '1.220'
$ws.Range("E34").Value = '  -1.33%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText "D35" '3.101'
$ws.Range("E35").Value = '  -3.40%  '

# Row 36
Set-CellText "D36" '0.7700'
$ws.Range("E36").Value = '  -0.13%  '

# Row 37
Set-CellText "D37" '0.02080'
$ws.Range("E37").Value = '  +1.17%  '

# Row 38
Set-CellText "D38" '2.559'
$ws.Range("E38").Value = '  -0.17%  '

# Row 39
Set-CellText "D39" '3.074'
$ws.Range("E39").Value = '  +1.99%  '

# Row 40
Set-CellText "D40" '0.5586'
$ws.Range("E40").Value = '  +0.83%  '

# Row 41
$ws.Range("E41").Value = '  -0.79%  '

# Row 42
Set-CellText "D42" '6.734'
$ws.Range("E42").Value = '  -3.61%  '

# Row 43
Set-CellText "D43" '116.60'
$ws.Range("E43").Value = '  +4.57%  '

# Row 44
Set-CellText "D44" '8.674'
$ws.Range("E44").Value = '  +1.47%  '

# Row 45
Set-CellText "D45" '0.1508'
$ws.Range("E45").Value = '  -1.04%  '

# Row 46
Set-CellText "D46" '0.4800'
$ws.Range("E46").Value = '  -0.09%  '

# Row 47
Set-CellText "D47" '10.41'
$ws.Range("E47").Value = '  -1.90%  '

# Row 48
Set-CellText "D48" '1.000'
$ws.Range("E48").Value = '  -0.16%  '

# Row 49
Set-CellText "D49" '1.595'
$ws.Range("E49").Value = '  -2.79%  '

# Row 50
Set-CellText "D50" '66.28'
$ws.Range("E50").Value = '  -1.69%  '

# Row 51
Set-CellText "D51" '37.10'
$ws.Range("E51").Value = '  +0.17%  '
